# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting refreshed counts from the data source (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 216
$ws1.Range("F11").Value = 507
$ws1.Range("F14").Value = 12617
$ws1.Range("F15").Value = 5204

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 31

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 216
$ws4.Range("F12").Value = 507
$ws4.Range("F15").Value = 12617
$ws4.Range("F16").Value = 31
$ws4.Range("F18").Value = 5204
